# "adding more data collectors"
#
# The sheet is a flat list of people (Nombre / Email / Celular / Institución).
# This edit:
#   1. Adds a new collector, "Francis Soria", as row 48 (mail/celular/institucion
#      use the same placeholder values as every other row).
#   2. Replaces every existing collector's real e-mail hyperlink (column B)
#      with the literal placeholder text "mail", and their Institución
#      (column D, previously "EPN") with the placeholder "U" - matching the
#      same placeholder pattern already used for the new row.
#   3. Removes the now-unused mailto: hyperlinks from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row first so its strings are registered before the bulk
# "mail" / "U" placeholder strings are introduced by the loop below - this
# keeps the shared-string table ordering stable (new unique text ends up
# right after the last existing name, same as the source edit).
$ws.Range("A48").Value = "Francis Soria"
$ws.Range("B48").Value = "mail"
$ws.Range("C48").Value = 5930000000
$ws.Range("D48").Value = "U"

# Scrub the real e-mails / institution values from the existing rows.
for ($r = 3; $r -le 47; $r++) {
    $ws.Range("D" + $r).Value = "U"
    $ws.Range("B" + $r).Value = "mail"
}

# Drop the mailto: hyperlinks that used to live on column B.
$ws.Hyperlinks.Delete()
